# Runsheet update: add exclusions
#  - rename the "keep" label to "include" in the exclude/include column (col I)
#  - turn on AutoFilter for the include/exclude column (I2:I36)
#  - select I2:I36 (matches the new filter range)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Find the "exclude" column (header in row 1) so this is robust to layout ---
$headerRange = $ws.Range("A1:I1")
$excludeCol = 0
for ($c = 1; $c -le $headerRange.Columns.Count; $c++) {
    $h = $ws.Cells.Item(1, $c).Value2
    if ($h -eq "exclude") {
        $excludeCol = $c
    }
}
if ($excludeCol -eq 0) {
    $excludeCol = 9  # fall back to column I
}

# --- Find the used extent of that column below the header ---
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$firstDataRow = 2
$lastDataRow = 1
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, $excludeCol).Value2
    if ($v -ne $null) {
        $lastDataRow = $r
    }
}
if ($lastDataRow -lt $firstDataRow) {
    $lastDataRow = $firstDataRow
}

# --- Replace every "keep" value with "include" in that column ---
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $cell = $ws.Cells.Item($r, $excludeCol)
    if ($cell.Value2 -eq "keep") {
        $cell.Value = "include"
    }
}

# --- Turn on AutoFilter over the data range of that column ---
$filterRange = $ws.Range($ws.Cells.Item($firstDataRow, $excludeCol), $ws.Cells.Item($lastDataRow, $excludeCol))
$filterRange.AutoFilter()

# --- Register the hidden _FilterDatabase defined name (sheet-scoped), like Excel does ---
$filterAddress = $ws.Name + "!" + $filterRange.Address($true, $true)
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $filterAddress)
$fdName.Visible = $false

# --- Update the saved selection to match the new filter range ---
$filterRange.Select()

Write-Host "Updated exclude/include column: rows $firstDataRow-$lastDataRow in column $excludeCol"
